# Generate Report for Handoff
# Updates the localization status report: files that were queued with "low"
# priority have now been handed off at a higher ("ht") priority, and their
# "Latest Handoff Datetime" is refreshed to the new handoff run's timestamp.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 4-7 correspond to the files that were still at
# "Ready for handoff" / low priority; bump priority and handoff time.
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($row in 4..7) {
    $wsZh.Cells.Item($row, 5).Value = "ht"
    $wsZh.Cells.Item($row, 8).Value = "2016-08-25 12:33:35"
}

# de-de sheet: same handoff run, different locale timestamp.
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($row in 4..7) {
    $wsDe.Cells.Item($row, 5).Value = "ht"
    $wsDe.Cells.Item($row, 8).Value = "2016-08-25 12:33:40"
}
